$d = $word.ActiveDocument

$range = $d.Content
$range.Find.Execute("{fechanota}", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)
